$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.234.51"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "2.235.84"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.17"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.627"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.69"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("E9").Value = "  -2.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.41"
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0960"
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.97"
$ws.Range("E12").Value = "  -2.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").Value = "2.572.60"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.37"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.840"
$ws.Range("E16").Value = "  -1.94%  "
$ws.Range("D17").Value = "2.231.15"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "42.104.67"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000105"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.23"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.88"
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.38"
$ws.Range("E22").Value = "  +7.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.40"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("E24").Value = "  -5.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.37"
$ws.Range("E26").Value = "  -2.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.62"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  -2.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.42"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.62"
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.64"
$ws.Range("E32").Value = "  -4.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0805"
$ws.Range("E33").Value = "  -1.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.88"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.125"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.110"
$ws.Range("E36").Value = "  -7.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.36"
$ws.Range("E37").Value = "  -4.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0305"
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.25"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.14"
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.72"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.05"
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.200"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.74"
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.70"
$ws.Range("E45").Value = "  -1.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.100"
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.13"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.35"
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.17"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.68"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").Value = "2.443.59"
$ws.Range("E51").Value = "  -0.31%  "
